$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 79 (Angeleno Primera/Segunda dated 44335).
# New weekly data (dated 44939) for Black Amber and Larry Ann varieties is
# inserted before the old rows, which get pushed down to rows 84:85 intact.
# Insert 6 new blank rows at 78:83 - old row 78 (Angeleno/Primera) shifts to
# row 84, and old row 79 (Angeleno/Segunda) shifts to row 85.
$ws.Rows("78:83").Insert()

# Row 78 - newly inserted row: Black Amber / Especial
$ws.Range("A78").Value = 11
$ws.Range("B78").Value = "Vega Monumental Concepción"
$ws.Range("C78").Value = "Bíobío"
$ws.Range("D78").Value = 44939
$ws.Range("E78").Value = 8
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100103
$ws.Range("H78").Value = "Frutos de hueso (carozo)"
$ws.Range("I78").Value = 100103002
$ws.Range("J78").Value = "Ciruela"
$ws.Range("K78").Value = "Black Amber"
$ws.Range("L78").Value = "Especial"
$ws.Range("M78").Value = 150
$ws.Range("N78").Value = 14000
$ws.Range("O78").Value = 14000
$ws.Range("P78").Value = 14000
$ws.Range("Q78").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R78").Value = "Región de O'Higgins"
$ws.Range("S78").Value = 933
$ws.Range("T78").Value = 15

# Row 79 - newly inserted row: Black Amber / Primera
$ws.Range("A79").Value = 11
$ws.Range("B79").Value = "Vega Monumental Concepción"
$ws.Range("C79").Value = "Bíobío"
$ws.Range("D79").Value = 44939
$ws.Range("E79").Value = 8
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value = 100103
$ws.Range("H79").Value = "Frutos de hueso (carozo)"
$ws.Range("I79").Value = 100103002
$ws.Range("J79").Value = "Ciruela"
$ws.Range("K79").Value = "Black Amber"
$ws.Range("L79").Value = "Primera"
$ws.Range("M79").Value = 120
$ws.Range("N79").Value = 12000
$ws.Range("O79").Value = 12000
$ws.Range("P79").Value = 12000
$ws.Range("Q79").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R79").Value = "Región de O'Higgins"
$ws.Range("S79").Value = 800
$ws.Range("T79").Value = 15

# Row 80 - newly inserted row: Black Amber / Segunda
$ws.Range("A80").Value = 11
$ws.Range("B80").Value = "Vega Monumental Concepción"
$ws.Range("C80").Value = "Bíobío"
$ws.Range("D80").Value = 44939
$ws.Range("E80").Value = 8
$ws.Range("F80").Value = "Fruta"
$ws.Range("G80").Value = 100103
$ws.Range("H80").Value = "Frutos de hueso (carozo)"
$ws.Range("I80").Value = 100103002
$ws.Range("J80").Value = "Ciruela"
$ws.Range("K80").Value = "Black Amber"
$ws.Range("L80").Value = "Segunda"
$ws.Range("M80").Value = 150
$ws.Range("N80").Value = 11000
$ws.Range("O80").Value = 11000
$ws.Range("P80").Value = 11000
$ws.Range("Q80").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R80").Value = "Región de O'Higgins"
$ws.Range("S80").Value = 733
$ws.Range("T80").Value = 15

# Row 81 - newly inserted row: Larry Ann / Especial
$ws.Range("A81").Value = 11
$ws.Range("B81").Value = "Vega Monumental Concepción"
$ws.Range("C81").Value = "Bíobío"
$ws.Range("D81").Value = 44939
$ws.Range("E81").Value = 8
$ws.Range("F81").Value = "Fruta"
$ws.Range("G81").Value = 100103
$ws.Range("H81").Value = "Frutos de hueso (carozo)"
$ws.Range("I81").Value = 100103002
$ws.Range("J81").Value = "Ciruela"
$ws.Range("K81").Value = "Larry Ann"
$ws.Range("L81").Value = "Especial"
$ws.Range("M81").Value = 120
$ws.Range("N81").Value = 14000
$ws.Range("O81").Value = 14000
$ws.Range("P81").Value = 14000
$ws.Range("Q81").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R81").Value = "Región de O'Higgins"
$ws.Range("S81").Value = 933
$ws.Range("T81").Value = 15

# Row 82 - newly inserted row: Larry Ann / Primera
$ws.Range("A82").Value = 11
$ws.Range("B82").Value = "Vega Monumental Concepción"
$ws.Range("C82").Value = "Bíobío"
$ws.Range("D82").Value = 44939
$ws.Range("E82").Value = 8
$ws.Range("F82").Value = "Fruta"
$ws.Range("G82").Value = 100103
$ws.Range("H82").Value = "Frutos de hueso (carozo)"
$ws.Range("I82").Value = 100103002
$ws.Range("J82").Value = "Ciruela"
$ws.Range("K82").Value = "Larry Ann"
$ws.Range("L82").Value = "Primera"
$ws.Range("M82").Value = 100
$ws.Range("N82").Value = 12000
$ws.Range("O82").Value = 12000
$ws.Range("P82").Value = 12000
$ws.Range("Q82").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R82").Value = "Región de O'Higgins"
$ws.Range("S82").Value = 800
$ws.Range("T82").Value = 15

# Row 83 - newly inserted row: Larry Ann / Segunda
$ws.Range("A83").Value = 11
$ws.Range("B83").Value = "Vega Monumental Concepción"
$ws.Range("C83").Value = "Bíobío"
$ws.Range("D83").Value = 44939
$ws.Range("E83").Value = 8
$ws.Range("F83").Value = "Fruta"
$ws.Range("G83").Value = 100103
$ws.Range("H83").Value = "Frutos de hueso (carozo)"
$ws.Range("I83").Value = 100103002
$ws.Range("J83").Value = "Ciruela"
$ws.Range("K83").Value = "Larry Ann"
$ws.Range("L83").Value = "Segunda"
$ws.Range("M83").Value = 120
$ws.Range("N83").Value = 11000
$ws.Range("O83").Value = 11000
$ws.Range("P83").Value = 11000
$ws.Range("Q83").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R83").Value = "Región de O'Higgins"
$ws.Range("S83").Value = 733
$ws.Range("T83").Value = 15

# Rows 84 (old row 78 content, Angeleno / Primera, untouched values) and 85
# (old row 79 content, Angeleno / Segunda, untouched values) were already
# shifted into place by the Insert() above and keep their original values.
